$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2268.4119
$ws.Range("J88").Value = 2431.9285
$ws.Range("L88").Value = 2431.9285
$ws.Range("N88").Value = -3243.9285
$ws.Range("H91").Value = 2268.4119
$ws.Range("J91").Value = 2431.9285
$ws.Range("L91").Value = 2431.9285
$ws.Range("N91").Value = -5239.9285

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6167.885
$ws.Range("I32").Value = 4464.7896
$ws.Range("K32").Value = 4464.7896
$ws.Range("M32").Value = -4177.7896
$ws.Range("H37").Value = 26257.25
$ws.Range("J37").Value = 31676.334
$ws.Range("L37").Value = 31676.334
$ws.Range("N37").Value = -32222.334
$ws.Range("H44").Value = 20049.334
$ws.Range("J44").Value = 20049.334
$ws.Range("L44").Value = 20049.334
$ws.Range("N44").Value = -21025.334
$ws.Range("H55").Value = 40841.5
$ws.Range("J55").Value = 44030
$ws.Range("L55").Value = 44030
$ws.Range("N55").Value = -44660
$ws.Range("H80").Value = 40405.555
$ws.Range("J80").Value = 38831.25
$ws.Range("L80").Value = 38831.25
$ws.Range("N80").Value = -40827.25
$ws.Range("H83").Value = 40405.555
$ws.Range("J83").Value = 38831.25
$ws.Range("L83").Value = 116493.75
$ws.Range("N83").Value = -126477.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8273.829
$ws.Range("I134").Value = 7211.857
$ws.Range("J134").Value = 9388.9
$ws.Range("K134").Value = 21635.571
$ws.Range("L134").Value = 28166.7
$ws.Range("M134").Value = -19100.571
$ws.Range("N134").Value = -33236.7

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 898.75
$ws.Range("J32").Value = 897.5
$ws.Range("L32").Value = 2692.5
$ws.Range("N32").Value = -3258.5
$ws.Range("H69").Value = 7177.88
$ws.Range("J69").Value = 9091
$ws.Range("L69").Value = 27273
$ws.Range("N69").Value = -28895
$ws.Range("H72").Value = 7177.88
$ws.Range("J72").Value = 9091
$ws.Range("L72").Value = 81819
$ws.Range("N72").Value = -89931

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1575
$ws.Range("I97").Value = 1109.5
$ws.Range("J97").Value = 2506
$ws.Range("K97").Value = 1109.5
$ws.Range("L97").Value = 2506
$ws.Range("M97").Value = -613.5
$ws.Range("N97").Value = -3498
$ws.Range("H122").Value = 3090.425
$ws.Range("I122").Value = 3255.32
$ws.Range("J122").Value = 2815.6
$ws.Range("K122").Value = 9765.960000000001
$ws.Range("L122").Value = 8446.799999999999
$ws.Range("M122").Value = -7315.960000000001
$ws.Range("N122").Value = -13346.8
$ws.Range("H126").Value = 3083.8
$ws.Range("I126").Value = 3450
$ws.Range("J126").Value = 2992.25
$ws.Range("K126").Value = 10350
$ws.Range("L126").Value = 8976.75
$ws.Range("M126").Value = -7880
$ws.Range("N126").Value = -13916.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 20839.166
$ws.Range("J24").Value = 20839.166
$ws.Range("L24").Value = 20839.166
$ws.Range("N24").Value = -21525.166
$ws.Range("H25").Value = 7483.7144
$ws.Range("I25").Value = 7064.5835
$ws.Range("K25").Value = 7064.5835
$ws.Range("M25").Value = -6834.5835
$ws.Range("H26").Value = 5494.5
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -705
$ws.Range("H31").Value = 3248.5789
$ws.Range("I31").Value = 964.0833
$ws.Range("J31").Value = 7164.857
$ws.Range("K31").Value = 964.0833
$ws.Range("L31").Value = 7164.857
$ws.Range("M31").Value = -716.0833
$ws.Range("N31").Value = -7660.857
$ws.Range("H34").Value = 10637
$ws.Range("I34").Value = 11250
$ws.Range("K34").Value = 11250
$ws.Range("M34").Value = -11078
$ws.Range("H43").Value = 587000
$ws.Range("I43").Value = 625000
$ws.Range("K43").Value = 625000
$ws.Range("M43").Value = -624807
$ws.Range("H132").Value = 37041700
$ws.Range("I132").Value = 55560556
$ws.Range("J132").Value = 3990
$ws.Range("K132").Value = 166681668
$ws.Range("L132").Value = 11970
$ws.Range("M132").Value = -166679138
$ws.Range("N132").Value = -17030

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 53198.6
$ws.Range("I37").Value = 41498.5
$ws.Range("K37").Value = 41498.5
$ws.Range("M37").Value = -41295.5
$ws.Range("H96").Value = 4374.25
$ws.Range("J96").Value = 3002
$ws.Range("L96").Value = 3002
$ws.Range("N96").Value = -5748
$ws.Range("H119").Value = 131099
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 131099
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 131099
$ws.Range("N119").Value = -140775
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 17356.125
$ws.Range("I122").Value = 17407
$ws.Range("J122").Value = 17000
$ws.Range("K122").Value = 52221
$ws.Range("L122").Value = 51000
$ws.Range("M122").Value = -49771
$ws.Range("N122").Value = -55900
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H124").Value = 48171.6
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 48171.6
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 48171.6
$ws.Range("N124").Value = -57991.6
$ws.Range("H125").Value = 49992.918
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 49992.918
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 49992.918
$ws.Range("N125").Value = -59832.918
$ws.Range("H126").Value = 3465.3442
$ws.Range("I126").Value = 4127.3516
$ws.Range("J126").Value = 2444.75
$ws.Range("K126").Value = 12382.0548
$ws.Range("L126").Value = 7334.25
$ws.Range("M126").Value = -9912.0548
$ws.Range("N126").Value = -12274.25
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 88000
$ws.Range("I129").Value = 88000
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 88000
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -83000
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 27921.447
$ws.Range("I132").Value = 22784.773
$ws.Range("J132").Value = 34984.375
$ws.Range("K132").Value = 68354.319
$ws.Range("L132").Value = 104953.125
$ws.Range("M132").Value = -65824.319
$ws.Range("N132").Value = -110013.125
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 1846.3077
$ws.Range("I136").Value = 1951.4286
$ws.Range("J136").Value = 1723.6666
$ws.Range("K136").Value = 5854.2858
$ws.Range("L136").Value = 5170.9998
$ws.Range("M136").Value = -3304.2858
$ws.Range("N136").Value = -10270.9998
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 49715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49715
$ws.Range("N139").Value = -59995
$ws.Range("H140").Value = 59999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 59999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
